$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 (TC00001): swap row1/row2 content to upper-case menu labels
# followed by the About/Support/Change Password/Logout row, and move
# the selection from D2 to D1.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A1").Value = "MENU1"
$ws1.Range("B1").Value = "MENU2"
$ws1.Range("C1").Value = "MENU3"
$ws1.Range("D1").Value = "MENU4"

$ws1.Range("A2").Value = "About"
$ws1.Range("B2").Value = "Support"
$ws1.Range("C2").Value = "Change Password"
$ws1.Range("D2").Value = "Logout"

$ws1.Range("D1").Select()

# ---------------------------------------------------------------------
# Sheet 2 (TC00002): SUPPORT_URL / CHANGE_PASSWORD_URL reference sheet
# with a live hyperlink in A2.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "TC00002"

$ws2.Range("A1").Value = "SUPPORT_URL"
$ws2.Range("B1").Value = "CHANGE_PASSWORD_URL"
$ws2.Range("A2").Value = "https://opensource-demo.orangehrmlive.com/web/index.php/help/support"
$ws2.Range("B2").Value = "https://opensource-demo.orangehrmlive.com/web/index.php/pim/updatePassword"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://opensource-demo.orangehrmlive.com/web/index.php/help/support")

$ws2.Columns.Item(1).ColumnWidth = 63.6
$ws2.Columns.Item(2).ColumnWidth = 70.9

$ws2.Range("B1").Select()

# ---------------------------------------------------------------------
# Sheet 3 (TC00003): side-menu headers + menu names, becomes the active
# tab.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "TC00003"

$ws3.Range("A1").Value = "SIDEMENU1"
$ws3.Range("B1").Value = "SIDEMENU2"
$ws3.Range("C1").Value = "SIDEMENU3"
$ws3.Range("D1").Value = "SIDEMENU4"
$ws3.Range("E1").Value = "SIDEMENU5"
$ws3.Range("F1").Value = "SIDEMENU6"
$ws3.Range("G1").Value = "SIDEMENU7"
$ws3.Range("H1").Value = "SIDEMENU8"
$ws3.Range("I1").Value = "SIDEMENU9"
$ws3.Range("J1").Value = "SIDEMENU10"
$ws3.Range("K1").Value = "SIDEMENU11"

$ws3.Range("A2").Value = "Admin"
$ws3.Range("B2").Value = "PIM"
$ws3.Range("C2").Value = "Leave"
$ws3.Range("D2").Value = "Time"
$ws3.Range("E2").Value = "Recruitment"
$ws3.Range("F2").Value = "My Info"
$ws3.Range("G2").Value = "Performance"
$ws3.Range("H2").Value = "Dashboard"
$ws3.Range("I2").Value = "Directory"
$ws3.Range("J2").Value = "Maintenance"
$ws3.Range("K2").Value = "Buzz"

$ws3.Columns.Item(10).ColumnWidth = 11.4
$ws3.Columns.Item(11).ColumnWidth = 11.4

$ws3.Range("K2").Select()
$ws3.Activate()
